$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.763.36"
$ws.Range("D3").Value = "2.278.00"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "251.44"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("D7").Value = "71.95"
$ws.Range("E7").Value = "  +8.26%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +12.31%  "
$ws.Range("D10").Value = "38.49"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").Value = "59.69"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").Value = "7.35"
$ws.Range("E13").Value = "  +6.98%  "
$ws.Range("D15").Value = "2.615.70"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "14.96"
$ws.Range("E16").Value = "  +4.19%  "
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("D18").Value = "2.270.47"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("D19").Value = "42.692.23"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("E20").Value = "  +6.94%  "
$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "236.31"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("D25").Value = "3.86"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").Value = "11.65"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").Value = "2.13"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").Value = "168.30"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "21.03"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("E33").Value = "  +11.20%  "
$ws.Range("E34").Value = "  +5.90%  "
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("D36").Value = "30.96"
$ws.Range("E36").Value = "  +26.32%  "
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("E38").Value = "  +16.60%  "
$ws.Range("D39").Value = "4.75"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "2.33"
$ws.Range("E41").Value = "  +5.06%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "13.31"
$ws.Range("E42").Value = "  +17.16%  "
$ws.Range("D43").Value = "5.86"
$ws.Range("E43").Value = "  +6.57%  "
$ws.Range("D44").Value = "0.211"
$ws.Range("E44").Value = "  +11.31%  "
$ws.Range("D45").Value = "9.16"
$ws.Range("E45").Value = "  +7.48%  "
$ws.Range("D46").Value = "4.97"
$ws.Range("E46").Value = "  -6.26%  "
$ws.Range("D47").Value = "61.57"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +3.98%  "
